$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.783183336257935
$ws.Range("B1").Value = 4.130315780639648
$ws.Range("C1").Value = 7.863772869110107
$ws.Range("D1").Value = 8.045058250427246
$ws.Range("E1").Value = 5.663721084594727
